# Update the date on the title line
$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-27 Saturday", "2024-07-28 Sunday"),
    @("289÷8=", "871÷6="),
    @("146÷7=", "707÷6="),
    @("992÷9=", "874÷2="),
    @("788÷7=", "450÷7="),
    @("595÷9=", "343÷5="),
    @("220÷3=", "401÷7="),
    @("943÷9=", "568÷2="),
    @("572÷3=", "548÷6="),
    @("137÷3=", "145÷2="),
    @("335÷5=", "198÷9="),
    @("666÷4=", "175÷5="),
    @("724÷4=", "317÷6="),
    @("153÷7=", "313÷7="),
    @("304÷5=", "890÷3="),
    @("326÷2=", "900÷3="),
    @("691÷4=", "894÷8="),
    @("364÷6=", "209÷5="),
    @("993÷7=", "518÷9="),
    @("537÷4=", "899÷5="),
    @("915÷6=", "194÷7="),
    @("114÷8=", "429÷8="),
    @("910÷8=", "610÷5="),
    @("639÷8=", "721÷9="),
    @("867÷9=", "577÷2="),
    @("211÷4=", "966÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
